$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.816.10"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "2.315.71"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "271.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "2.658.50"
$ws.Range("E14").Value = "  +4.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.873"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.73%  "
$ws.Range("D17").Value = "2.317.06"
$ws.Range("E17").Value = "  +4.20%  "
$ws.Range("D18").Value = "43.758.80"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E19").Value = "  +4.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.127"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("E36").Value = "  +3.48%  "
$ws.Range("E37").Value = "  -3.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("E39").Value = "  -6.47%  "
$ws.Range("E40").Value = "  +12.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.14%  "
$ws.Range("E42").Value = "  +18.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "62.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +4.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.190"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +16.27%  "
$ws.Range("D51").Value = "2.545.25"
$ws.Range("E51").Value = "  +4.13%  "
